$d = $word.ActiveDocument

# Remove the hidden _GoBack bookmark from its current location (end of the
# "DrawBoard ... square class." paragraph); it will be re-added at the end
# of the newly appended content below.
try {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
} catch {
    # no existing _GoBack bookmark - nothing to remove
}

# Collapsed range sitting right at the very end of the document's content
# (after the last character, before the implicit end-of-story mark) so that
# InsertXML appends new paragraphs rather than splitting/duplicating the
# existing final paragraph.
$insertAt = $d.Content.End
$r = $d.Range($insertAt, $insertAt)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListBullet"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="0"/>
              </w:numPr>
              <w:ind w:left="360" w:hanging="360"/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListBullet"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="0"/>
              </w:numPr>
              <w:ind w:left="360" w:hanging="360"/>
              <w:rPr>
                <w:b/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
              </w:rPr>
              <w:t>09/12/2016</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListBullet"/>
            </w:pPr>
            <w:r>
              <w:t>Used pointers to solve the issue of linking the squares and their corresponding dices, if any.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListBullet"/>
            </w:pPr>
            <w:r>
              <w:t>Some modification within the square and dice classes to</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> integrate pointers to setup a</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> link</w:t>
            </w:r>
            <w:r>
              <w:t>.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListBullet"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Completed </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>DrawBoard</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> and </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>UpdateBoard</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> function in the Board class. Also, integrated a multidimensional string array that will be utilized later for serialization.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListBullet"/>
            </w:pPr>
            <w:r>
              <w:t>The functions within Board class that print results to the console should be separated in a view class.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListBullet"/>
            </w:pPr>
            <w:r>
              <w:t>Having issues with coming up with a proper way to access the Board class objects from other classes without compromising on Data Encapsulation.</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$countBefore = $d.Paragraphs.Count
$r.InsertXML($xml) | Out-Null
$countAfter = $d.Paragraphs.Count
$addedCount = $countAfter - $countBefore

# The literal <w:ind w:left="360" w:hanging="360"/> in the payload above gets
# normalized away by the save pipeline because it is identical to the
# indentation the ListBullet style's linked numbering level (numId 1, ilvl 0)
# would already contribute - even though these two paragraphs explicitly
# opt out of that numbering via numId=0. Re-apply the same indentation as
# direct paragraph formatting so it survives serialization.
$firstNewIndex = $countBefore + 1
$blankPara = $d.Paragraphs($firstNewIndex)
$blankPara.LeftIndent = 18
$blankPara.FirstLineIndent = -18

$datePara = $d.Paragraphs($firstNewIndex + 1)
$datePara.LeftIndent = 18
$datePara.FirstLineIndent = -18
